# New changes - 6/30/2020
# Splits the single "Sheet1" mapping sheet into two sheets:
#   - Bank_Report        (A1:C3 - FolderPath/FileType/Company bank-report info)
#   - Remittance_Report  (A1:E3 - the original 5-column remittance mapping, re-pointed
#                          at the new UiPath folder layout)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the existing sheet so we end up with two sheets that both carry
#    the original formatting/styles, then rename them.
# ---------------------------------------------------------------------------
$bank = $wb.Worksheets.Item(1)
$bank.Copy($null, $bank)
$remit = $wb.Worksheets.Item(2)

$bank.Name = "Bank_Report"
$remit.Name = "Remittance_Report"

# ---------------------------------------------------------------------------
# 2. Remittance_Report keeps the 5-column shape of the original sheet - just
#    refresh the cell values.
# ---------------------------------------------------------------------------
$remit.Cells.Item(1,1).Value = "Client Name"
$remit.Cells.Item(1,2).Value = "FileType"
$remit.Cells.Item(1,3).Value = "FolderPath"
$remit.Cells.Item(1,4).Value = "Company name"
$remit.Cells.Item(1,5).Value = "Biller"

$remit.Cells.Item(2,1).Value = "Element"
$remit.Cells.Item(2,2).Value = "Excel"
$remit.Cells.Item(2,3).Value = "C:\Users\Hp\Documents\UiPath\AR2.0\Data\Input\RemittanceReport\"
$remit.Cells.Item(2,4).Value = "ELEMENT FLEET"
$remit.Cells.Item(2,5).Value = "Freddy"

$remit.Cells.Item(3,1).Value = "Element_Test"
$remit.Cells.Item(3,2).Value = "Excel"
$remit.Cells.Item(3,3).Value = "C:\Users\Hp\Documents\UiPath\AR2.0\Data\Input\RemittanceReport\"
$remit.Cells.Item(3,4).Value = "ELEMENT FLEET"
$remit.Cells.Item(3,5).Value = "Freddy"

# Normalise the data-row styling so every body cell shares the same
# plain-bordered look (copy the already-correct format from A2 over the rest
# of the data rows).
$remit.Range("A2").Copy() | Out-Null
$remit.Range("B2:E3").PasteSpecial(-4122) | Out-Null
$remit.Range("A2").Copy() | Out-Null
$remit.Range("A3").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3. Bank_Report shrinks down to 3 columns with brand new content.
# ---------------------------------------------------------------------------
$bank.Cells.Item(1,1).Value = "Bank Report"
$bank.Cells.Item(1,2).Value = "FileType"
$bank.Cells.Item(1,3).Value = "FolderPath"

$bank.Cells.Item(2,1).Value = "Bridge Bank"
$bank.Cells.Item(2,2).Value = "Excel"
$bank.Cells.Item(2,3).Value = "C:\Users\Hp\Documents\UiPath\AR2.0\Data\Input\BankReport\"

$bank.Cells.Item(3,1).Value = "PNC Bank"
$bank.Cells.Item(3,2).Value = "Excel"
$bank.Cells.Item(3,3).Value = "C:\Users\Hp\Documents\UiPath\AR2.0\Data\Input\BankReport\"

# Drop the no-longer-used D/E columns entirely.
$bank.Columns("D:E").Delete() | Out-Null

# Normalise styling on the remaining body rows the same way as above.
$bank.Range("A2").Copy() | Out-Null
$bank.Range("B2:C3").PasteSpecial(-4122) | Out-Null
$bank.Range("A2").Copy() | Out-Null
$bank.Range("A3").PasteSpecial(-4122) | Out-Null

$bank.Columns("A:C").AutoFit() | Out-Null
$remit.Columns("A:E").AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore the view state: Bank_Report is the active tab with A3 selected.
# ---------------------------------------------------------------------------
$bank.Activate()
$bank.Range("A3").Select()

Write-Output "done"
